$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.247.04"
$ws.Range("E2").Value = "  -2.14%  "

$ws.Range("D3").Value = "2.428.88"
$ws.Range("E3").Value = "  -1.81%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.25%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "565.51"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.01%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.00"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.99%  "

$ws.Range("E7").Value = "  +0.21%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.527"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.65%  "

$ws.Range("D9").Value = "2.424.21"
$ws.Range("E9").Value = "  -2.30%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.107"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -6.21%  "

$ws.Range("E11").Value = "  +0.61%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.17"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.77%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.351"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.58%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.45"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.11%  "

$ws.Range("E15").Value = "  -6.37%  "

$ws.Range("D16").Value = "2.869.19"
$ws.Range("E16").Value = "  -1.24%  "

$ws.Range("D17").Value = "62.080.49"
$ws.Range("E17").Value = "  -1.98%  "

$ws.Range("D18").Value = "2.419.34"
$ws.Range("E18").Value = "  -2.35%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.00"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.91%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.09"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.56%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "323.77"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.53%  "

$ws.Range("E22").Value = "  -3.16%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.42%  "

$ws.Range("E24").Value = "  +0.41%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "65.02"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.67%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "616.71"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.39%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.02"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.21%  "

$ws.Range("D28").Value = "0.0₃0954"
$ws.Range("E28").Value = "  -9.87%  "

$ws.Range("D29").Value = "2.545.38"
$ws.Range("E29").Value = "  -1.94%  "

$ws.Range("E30").Value = "  +0.43%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.43"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -6.25%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.99"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.25%  "

$ws.Range("E33").Value = "  -3.95%  "

$ws.Range("E34").Value = "  -8.28%  "

$ws.Range("E35").Value = "  -3.92%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.21%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.44"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -7.43%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.374"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.17%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.47"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.83%  "

$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.20"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -6.57%  "

$ws.Range("B41").Value = "Monero"
$ws.Range("C41").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "145.81"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.49%  "

$ws.Range("B42").Value = "OKB"
$ws.Range("C42").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "42.77"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.26%  "

$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.71"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -7.85%  "

$ws.Range("E44").Value = "  +0.01%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.45"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -8.95%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "143.75"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.79%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.68"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.65%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0521"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.68%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "19.99"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.76%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.591"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.21%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0228"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.36%  "
